$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record for "Pina / Caramelo" was added to the series.
# It belongs between the current rows 133 and 134 (chronologically it sits
# right before the old row 134), so insert a fresh row at 134 and push
# everything from the old row 134 onward down by one.
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with the new record's data. All of
# the "constant" columns for this sub-series (market, region, product
# hierarchy, origin, etc.) match every neighboring row.
$ws.Cells.Item(134, 1).Value = 5
$ws.Cells.Item(134, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(134, 3).Value = "Maule"
$ws.Cells.Item(134, 4).Value = 44461
$ws.Cells.Item(134, 5).Value = 7
$ws.Cells.Item(134, 6).Value = "Fruta"
$ws.Cells.Item(134, 7).Value = 100108
$ws.Cells.Item(134, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(134, 9).Value = 100108005
$ws.Cells.Item(134, 10).Value = "Piña"
$ws.Cells.Item(134, 11).Value = "Caramelo"
$ws.Cells.Item(134, 12).Value = "Segunda"
$ws.Cells.Item(134, 13).Value = 180
$ws.Cells.Item(134, 14).Value = 19000
$ws.Cells.Item(134, 15).Value = 19000
$ws.Cells.Item(134, 16).Value = 19000
$ws.Cells.Item(134, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(134, 18).Value = "Ecuador"
$ws.Cells.Item(134, 19).Value = 1357
$ws.Cells.Item(134, 20).Value = 14
